$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.746.52'
$ws.Range('E2').Value = '  +5.19%  '
$ws.Range('D3').Value = '2.653.22'
$ws.Range('E3').Value = '  +5.76%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.04'
$ws.Range('E5').Value = '  +2.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.21'
$ws.Range('E6').Value = '  +3.19%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +2.78%  '
$ws.Range('E9').Value = '  +15.51%  '
$ws.Range('D10').Value = '2.653.48'
$ws.Range('E10').Value = '  +5.85%  '
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('E12').Value = '  +4.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.10'
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000192'
$ws.Range('E14').Value = '  +10.34%  '
$ws.Range('D15').Value = '3.147.15'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.98'
$ws.Range('E16').Value = '  +5.07%  '
$ws.Range('D17').Value = '72.511.74'
$ws.Range('E17').Value = '  +5.04%  '
$ws.Range('D18').Value = '2.651.79'
$ws.Range('E18').Value = '  +5.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '385.48'
$ws.Range('E19').Value = '  +6.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.61'
$ws.Range('E20').Value = '  +6.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.94'
$ws.Range('E21').Value = '  +5.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.23'
$ws.Range('E22').Value = '  +5.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.04'
$ws.Range('E23').Value = '  +22.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.35'
$ws.Range('E24').Value = '  +4.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.45'
$ws.Range('E25').Value = '  +6.68%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.00'
$ws.Range('E27').Value = '  +11.61%  '
$ws.Range('D28').Value = '2.791.60'
$ws.Range('E28').Value = '  +6.09%  '
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('D30').Value = '0.0₃0975'
$ws.Range('E30').Value = '  +10.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '543.94'
$ws.Range('E31').Value = '  +7.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.10'
$ws.Range('E32').Value = '  +5.02%  '
$ws.Range('E33').Value = '  +11.07%  '
$ws.Range('E34').Value = '  +4.76%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '165.39'
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.44'
$ws.Range('E37').Value = '  +4.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.114'
$ws.Range('E38').Value = '  -3.04%  '
$ws.Range('E39').Value = '  +8.99%  '
$ws.Range('E40').Value = '  +2.36%  '
$ws.Range('E41').Value = '  +10.28%  '
$ws.Range('E42').Value = '  +7.94%  '
$ws.Range('E43').Value = '  +15.45%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('E45').Value = '  +5.40%  '
$ws.Range('E46').Value = '  +2.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '152.27'
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('E48').Value = '  +4.54%  '
$ws.Range('E49').Value = '  +6.68%  '
$ws.Range('E50').Value = '  +10.24%  '
$ws.Range('D51').Value = '0.0₆0270'
$ws.Range('E51').Value = '  +11.62%  '
